# Revision History - Figma.xlsx
# Add a new revision-history row (row 12): 06/11/2022, version 0.4,
# "Aggiunta StateChart", author SDS.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C12: version string "0.4" -------------------------------------------
# Force text storage (so "0.4" doesn't get stored as the number 0.4),
# then restore the original cell style (center aligned, same as the rest
# of the VERSIONE column) by pasting formats from the row above.
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "0.4"
[void]$ws.Range("C11").Copy()
[void]$ws.Range("C12").PasteSpecial(-4122)

# --- D12: description -----------------------------------------------------
$ws.Range("D12").Value = "Aggiunta StateChart"

# --- E12: author ------------------------------------------------------------
$ws.Range("E12").Value = "SDS"

# --- B12: date (06 Nov 2022 -> serial 44871) -------------------------------
# Copy the date-formatted style from the row above first so the new cell
# matches the existing date column formatting, then set the value.
[void]$ws.Range("B11").Copy()
[void]$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value = 44871

# --- Update the active selection to match the reviewed range --------------
[void]$ws.Range("B3:B12").Select()

# --- Column widths: the VERSIONE column (B) is widened and no longer
# relies on auto best-fit; the others keep their existing best-fit widths.
$ws.Columns.Item(2).ColumnWidth = 11.917

Write-Host "Revision row added"
